$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the "lab_sample_id" row (row 3 in the original layout)
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).Delete()
# Layout now:
#  1 header
#  2 participant_id
#  3 platekey
#  4 delivery_id
#  5 delivery_date
#  6 genome_build
#  7 type
#  8 file_path
#  9 filename
# 10 file_sub_type
# 11 file_type
# 12 delivery_version

# ---------------------------------------------------------------------------
# 2. Insert two new rows after "platekey" (row 3) for the new fields
#    referral_id / associated_interpretation_request_id
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()
$ws.Rows.Item(4).Insert()

# Copy the formatting of the platekey row (row 3) onto the two fresh rows
$ws.Range("A3:D3").Copy()
$ws.Range("A4:D5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(4,1).Value = "genome_file_paths_and_types"
$ws.Cells.Item(4,2).Value = "referral_id"
$ws.Cells.Item(4,3).Value = "varchar"
$ws.Cells.Item(4,4).Value = ""

$ws.Cells.Item(5,1).Value = "genome_file_paths_and_types"
$ws.Cells.Item(5,2).Value = "associated_interpretation_request_id"
$ws.Cells.Item(5,3).Value = "varchar"
$ws.Cells.Item(5,4).Value = ""

# Layout now:
#  1 header
#  2 participant_id
#  3 platekey
#  4 referral_id (new)
#  5 associated_interpretation_request_id (new)
#  6 delivery_id
#  7 delivery_date
#  8 genome_build
#  9 type
# 10 file_path
# 11 filename
# 12 file_sub_type
# 13 file_type
# 14 delivery_version

# ---------------------------------------------------------------------------
# 3. Rename field "file_path" -> "path" (row 10)
# ---------------------------------------------------------------------------
$ws.Cells.Item(10,2).Value = "path"

# ---------------------------------------------------------------------------
# 4. Remove the "file_type" row (row 13)
# ---------------------------------------------------------------------------
$ws.Rows.Item(13).Delete()

# Layout now:
#  1 header
#  2 participant_id
#  3 platekey
#  4 referral_id
#  5 associated_interpretation_request_id
#  6 delivery_id
#  7 delivery_date
#  8 genome_build
#  9 type
# 10 path
# 11 filename
# 12 file_sub_type
# 13 delivery_version

# ---------------------------------------------------------------------------
# 5. Append two new rows: software_version / delivery_type
# ---------------------------------------------------------------------------
$ws.Range("A13:D13").Copy()
$ws.Range("A14:D15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Give the two new trailing rows their own border treatment (left/right thin
# only, no top/bottom) - matches the new style introduced for these rows.
$ws.Range("A14:D15").Borders.Item(8).LineStyle = -4142
$ws.Range("A14:D15").Borders.Item(9).LineStyle = -4142

$ws.Cells.Item(14,1).Value = "genome_file_paths_and_types"
$ws.Cells.Item(14,2).Value = "software_version"
$ws.Cells.Item(14,3).Value = "varchar"
$ws.Cells.Item(14,4).Value = "newly added field"

$ws.Cells.Item(15,1).Value = "genome_file_paths_and_types"
$ws.Cells.Item(15,2).Value = "delivery_type"
$ws.Cells.Item(15,3).Value = "varchar"
$ws.Cells.Item(15,4).Value = "newly added field"

# ---------------------------------------------------------------------------
# 6. Row heights, matching the refreshed layout
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 31
$ws.Rows.Item(2).RowHeight = 155
$ws.Rows.Item(3).RowHeight = 170.5
$ws.Rows.Item(4).RowHeight = 77.5
$ws.Rows.Item(5).RowHeight = 77.5
$ws.Rows.Item(6).RowHeight = 139.5
$ws.Rows.Item(7).RowHeight = 108.5
$ws.Rows.Item(8).RowHeight = 108.5
$ws.Rows.Item(9).RowHeight = 124
$ws.Rows.Item(10).RowHeight = 77.5
$ws.Rows.Item(11).RowHeight = 77.5
$ws.Rows.Item(12).RowHeight = 108.5
$ws.Rows.Item(13).RowHeight = 77.5
$ws.Rows.Item(14).RowHeight = 77.5
$ws.Rows.Item(15).RowHeight = 77.5

# ---------------------------------------------------------------------------
# 7. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 20.54296875
$ws.Columns.Item(3).ColumnWidth = 30.1796875

# ---------------------------------------------------------------------------
# 8. View selection - move selection/active cell to G13, drop the frozen
#    top-left scroll position
# ---------------------------------------------------------------------------
$ws.Range("G13").Select()
